$wb = $excel.ActiveWorkbook

# Update "Created at" timestamp on the Input Values sheet
$inputSheet = $wb.Worksheets.Item("Input Values")
$inputSheet.Range("C2").Value = "2022-12-12 18:30:02.933102"

# Update the ApPredict version information sheet
$verSheet = $wb.Worksheets.Item("ApPredict version information")

# Replace the ApPredict arguments value (B15) with multi-line text
$argsText = @'
--pacing-freq 1 
--pacing-max-time 5 
--plasma-conc-high 100 
--plasma-conc-low 0 
--plasma-conc-count 4 
--plasma-conc-logscale true 
--model 1
'@
$verSheet.Range("B15").Value = $argsText

# Add a new row with Python packages version info
$verSheet.Range("A16").Value = "Python packages versions for chaste_codegen"

$pythonPackagesInfo = @'
Python 3.7.3
cellmlmanip==0.3.5
chaste-codegen==0.8.0
decorator==4.4.2
importlib-metadata==4.13.0
isodate==0.6.1
Jinja2==2.11.3
lxml==4.9.1
MarkupSafe==1.1.1
mpmath==1.2.1
networkx==2.6.3
packaging==20.9
Pint==0.18
pkg_resources==0.0.0
py==1.11.0
pyparsing==2.4.7
rdflib==5.0.0
six==1.16.0
sympy==1.10.1
typing_extensions==4.4.0
zipp==1.2.0
$CHASTE_TEST_OUTPUT is currently set to  /home/appredict/apps/ApPredict/testoutput.
Copyright (c) 2005-2021, University of Oxford.
All rights reserved.
University of Oxford means the Chancellor, Masters and Scholars of the
University of Oxford, having an administrative office at Wellington
Square, Oxford OX1 2JD, UK.
Redistribution and use in source and binary forms, with or without
modification, are permitted provided that the following conditions are met:
 * Redistributions of source code must retain the above copyright notice,
   this list of conditions and the following disclaimer.
 * Redistributions in binary form must reproduce the above copyright notice,
   this list of conditions and the following disclaimer in the documentation
   and/or other materials provided with the distribution.
 * Neither the name of the University of Oxford nor the names of its
   contributors may be used to endorse or promote products derived from this
   software without specific prior written permission.
THIS SOFTWARE IS PROVIDED BY THE COPYRIGHT HOLDERS AND CONTRIBUTORS "AS IS"
AND ANY EXPRESS OR IMPLIED WARRANTIES, INCLUDING, BUT NOT LIMITED TO, THE
IMPLIED WARRANTIES OF MERCHANTABILITY AND FITNESS FOR A PARTICULAR PURPOSE
ARE DISCLAIMED. IN NO EVENT SHALL THE COPYRIGHT HOLDER OR CONTRIBUTORS BE
LIABLE FOR ANY DIRECT, INDIRECT, INCIDENTAL, SPECIAL, EXEMPLARY, OR
CONSEQUENTIAL DAMAGES (INCLUDING, BUT NOT LIMITED TO, PROCUREMENT OF SUBSTITUTE
GOODS OR SERVICES; LOSS OF USE, DATA, OR PROFITS; OR BUSINESS INTERRUPTION)
HOWEVER CAUSED AND ON ANY THEORY OF LIABILITY, WHETHER IN CONTRACT, STRICT
LIABILITY, OR TORT (INCLUDING NEGLIGENCE OR OTHERWISE) ARISING IN ANY WAY OUT
OF THE USE OF THIS SOFTWARE, EVEN IF ADVISED OF THE POSSIBILITY OF SUCH DAMAGE.
This version of Chaste was compiled on:
Mon, 07 Nov 2022 17:49:14 +0000 by Linux d09b088bdc9f 4.15.0-161-generic #169-Ubuntu SMP Fri Oct 15 13:41:54 UTC 2021 x86_64 (uname)
from revision number 682dce0 with build type GccOpt, shared libraries.
ApPredict is based on commit 37cc5a6. But it HAS BEEN MODIFIED from that commit!
<ChasteBuildInfo>
	<ProvenanceInfo>
		<VersionString>2019.1.682dce0</VersionString> <!-- build specific -->
		<IsWorkingCopyModified>1</IsWorkingCopyModified>
		<BuildInformation>GccOpt, shared libraries</BuildInformation>
		<BuildTime>Mon, 07 Nov 2022 17:49:14 +0000</BuildTime>
		<CurrentTime>Thu, 10 Nov 2022 14:48:20 +0000</CurrentTime>
		<BuilderUnameInfo>Linux d09b088bdc9f 4.15.0-161-generic #169-Ubuntu SMP Fri Oct 15 13:41:54 UTC 2021 x86_64</BuilderUnameInfo>
		<Projects>
			<Project>
				<Name>ApPredict</Name>
				<Version>37cc5a6</Version>
				<Modified>True</Modified>
			</Project>
		</Projects>
	</ProvenanceInfo>
	<Compiler>
		<NameAndVersion>gcc, version b'9.3.0'</NameAndVersion>
		<Flags>-O3 -std=c++14</Flags>
	</Compiler>
	<Libraries>
		<CompiledIn>
			<PETSc>3.12.4</PETSc>
			<Boost>1.65.1</Boost>
			<HDF5>1.8.16</HDF5>
			<Parmetis>4.0.3</Parmetis>
		</CompiledIn>
		<Binaries>
			<XSD>4.0.0</XSD>
		</Binaries>
		<Optional>
			<SUNDIALS>2.5.0</SUNDIALS><!-- includes Cvode of a different version number -->
			<VTK>no</VTK>
			<Xerces>3.2.0</Xerces>
		</Optional>
	</Libraries>
</ChasteBuildInfo>
ApPredict args : --pacing-freq 1 --pacing-max-time 5 --plasma-conc-high 100 --plasma-conc-low 0 --plasma-conc-count 4 --plasma-conc-logscale true --model 1
HTTP Request : {"pacingFrequency":1,"pacingMaxTime":5,"plasmaMinimum":0,"plasmaMaximum":100,"plasmaIntermediatePointCount":"4","plasmaIntermediatePointLogScale":true,"modelId":"1"}
~/apps/app-manager/run/e716876a-4697-44b9-9267-4c6db505a40e ~/apps/app-manager

'@

$verSheet.Range("B16").Value = $pythonPackagesInfo
